$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 7693
$ws.Range("I80").Value = 987.3333
$ws.Range("J80").Value = 16314.571
$ws.Range("K80").Value = 2961.9999
$ws.Range("L80").Value = 48943.713
$ws.Range("M80").Value = -1963.9999
$ws.Range("N80").Value = -50939.713

$ws.Range("H83").Value = 7693
$ws.Range("I83").Value = 987.3333
$ws.Range("J83").Value = 16314.571
$ws.Range("K83").Value = 8885.9997
$ws.Range("L83").Value = 146831.139
$ws.Range("M83").Value = -3893.9997
$ws.Range("N83").Value = -156815.139

$ws.Range("H109").Value = 111899.164
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 111899.164
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 111899.164
$ws.Range("N109").Value = -114673.164

$ws.Range("H112").Value = 2416.3625
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 2468.0642
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 7404.192599999999
$ws.Range("M112").Value = -92
$ws.Range("N112").Value = -9620.192599999998

$ws.Range("H115").Value = 2000
$ws.Range("I115").Value = 2000
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 6000
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -4433
$ws.Range("N115").Value = $null

$ws.Range("H118").Value = 494.13333
$ws.Range("I118").Value = 379.3846
$ws.Range("J118").Value = 1240
$ws.Range("K118").Value = 1138.1538
$ws.Range("L118").Value = 3720
$ws.Range("M118").Value = 518.8462
$ws.Range("N118").Value = -7034

$ws.Range("H123").Value = 62067.406
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 62067.406
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 62067.406
$ws.Range("N123").Value = -71867.406

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null

$ws.Range("H126").Value = 45000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 45000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 45000
$ws.Range("N126").Value = -54880

$ws.Range("H127").Value = 1102.0435
$ws.Range("I127").Value = 665.6667
$ws.Range("J127").Value = 1116.7528
$ws.Range("K127").Value = 1997.0001
$ws.Range("L127").Value = 3350.2584
$ws.Range("M127").Value = 2962.9999
$ws.Range("N127").Value = -13270.2584

$ws.Range("H129").Value = 1010.5714
$ws.Range("I129").Value = 700
$ws.Range("J129").Value = 1017.0417
$ws.Range("K129").Value = 2100
$ws.Range("L129").Value = 3051.1251
$ws.Range("M129").Value = 2900
$ws.Range("N129").Value = -13051.1251

$ws.Range("H141").Value = 5194.971
$ws.Range("I141").Value = 2337.1333
$ws.Range("J141").Value = 22342
$ws.Range("K141").Value = 7011.3999
$ws.Range("L141").Value = 67026
$ws.Range("M141").Value = -1831.3999
$ws.Range("N141").Value = -77386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1139.125
$ws.Range("I61").Value = 1124.6957
$ws.Range("J61").Value = 1176
$ws.Range("K61").Value = 1124.6957
$ws.Range("L61").Value = 1176
$ws.Range("M61").Value = -912.6957
$ws.Range("N61").Value = -1600

$ws.Range("H88").Value = 2227.4285
$ws.Range("I88").Value = 1765.3334
$ws.Range("J88").Value = 2574
$ws.Range("K88").Value = 1765.3334
$ws.Range("L88").Value = 2574
$ws.Range("M88").Value = -1359.3334
$ws.Range("N88").Value = -3386

$ws.Range("H91").Value = 2227.4285
$ws.Range("I91").Value = 1765.3334
$ws.Range("J91").Value = 2574
$ws.Range("K91").Value = 1765.3334
$ws.Range("L91").Value = 2574
$ws.Range("M91").Value = -361.3334
$ws.Range("N91").Value = -5382

$ws.Range("H132").Value = 1941.2
$ws.Range("I132").Value = 1354.381
$ws.Range("J132").Value = 2821.4285
$ws.Range("K132").Value = 4063.143
$ws.Range("L132").Value = 8464.2855
$ws.Range("M132").Value = -1533.143
$ws.Range("N132").Value = -13524.2855

$ws.Range("H136").Value = 1139.125
$ws.Range("I136").Value = 1124.6957
$ws.Range("J136").Value = 1176
$ws.Range("K136").Value = 3374.0871
$ws.Range("L136").Value = 3528
$ws.Range("M136").Value = -824.0870999999997
$ws.Range("N136").Value = -8628

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 41394.332
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 41394.332
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 41394.332
$ws.Range("N108").Value = -49074.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 37000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 37000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 37000
$ws.Range("N88").Value = -37812

$ws.Range("H91").Value = 37000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 37000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 37000
$ws.Range("N91").Value = -39808

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = $null

$ws.Range("H134").Value = 1510.9574
$ws.Range("I134").Value = 1226.2354
$ws.Range("J134").Value = 2255.6155
$ws.Range("K134").Value = 3678.7062
$ws.Range("L134").Value = 6766.8465
$ws.Range("M134").Value = -1143.7062
$ws.Range("N134").Value = -11836.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 865.75
$ws.Range("I33").Value = 336.8
$ws.Range("J33").Value = 1747.3334
$ws.Range("K33").Value = 2020.8
$ws.Range("L33").Value = 10484.0004
$ws.Range("M33").Value = -1737.8
$ws.Range("N33").Value = -11050.0004

$ws.Range("H35").Value = 3599.75
$ws.Range("I35").Value = 900
$ws.Range("J35").Value = 4499.6665
$ws.Range("K35").Value = 2700
$ws.Range("L35").Value = 13498.9995
$ws.Range("M35").Value = -2412
$ws.Range("N35").Value = -14074.9995

$ws.Range("H137").Value = 23812632
$ws.Range("I137").Value = 2238.3333
$ws.Range("J137").Value = 41670428
$ws.Range("K137").Value = 6714.999899999999
$ws.Range("L137").Value = 125011284
$ws.Range("M137").Value = -1614.999899999999
$ws.Range("N137").Value = -125021484

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 26333.334
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 26333.334
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 26333.334
$ws.Range("N114").Value = -35011.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4586.095
$ws.Range("I7").Value = 3944.7778
$ws.Range("J7").Value = 5067.0835
$ws.Range("K7").Value = 3944.7778
$ws.Range("L7").Value = 5067.0835
$ws.Range("M7").Value = -3832.7778
$ws.Range("N7").Value = -5291.0835

$ws.Range("H126").Value = 4586.095
$ws.Range("I126").Value = 3944.7778
$ws.Range("J126").Value = 5067.0835
$ws.Range("K126").Value = 11834.3334
$ws.Range("L126").Value = 15201.2505
$ws.Range("M126").Value = -9364.3334
$ws.Range("N126").Value = -20141.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = $null
$ws.Range("N51").Value = $null
